$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.011.36"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.257.24"
$ws.Range("E3").Value = "  -3.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.22%  "
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.73%  "
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.16"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -9.27%  "
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "2.608.63"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.42%  "
$ws.Range("D17").Value = "2.266.93"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.773"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").Value = "42.023.51"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0692"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.79"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.87"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0986"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.64%  "
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("E42").Value = "  -8.89%  "
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "1.943.70"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.88%  "
$ws.Range("E48").Value = "  -5.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").Value = "2.482.79"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.50%  "
